$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.406.78'
$ws.Range("E2").Value = '  -1.80%  '

$ws.Range("D3").Value = '2.455.05'
$ws.Range("E3").Value = '  -1.93%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '560.98'
$ws.Range("E5").Value = '  -2.65%  '

$ws.Range("D6").Value = '163.59'
$ws.Range("E6").Value = '  -1.95%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = '0.504'
$ws.Range("E8").Value = '  -1.74%  '

$ws.Range("D9").Value = '2.455.21'
$ws.Range("E9").Value = '  -1.86%  '

$ws.Range("D10").Value = '0.151'
$ws.Range("E10").Value = '  -5.54%  '

$ws.Range("E11").Value = '  -1.98%  '

$ws.Range("E12").Value = '  -5.17%  '

$ws.Range("E13").Value = '  -2.35%  '

$ws.Range("D14").Value = '2.914.12'
$ws.Range("E14").Value = '  -1.56%  '

$ws.Range("D15").Value = '68.390.28'
$ws.Range("E15").Value = '  -1.60%  '

$ws.Range("D16").Value = '0.0000170'
$ws.Range("E16").Value = '  -3.63%  '

$ws.Range("D17").Value = '23.38'
$ws.Range("E17").Value = '  -5.53%  '

$ws.Range("D18").Value = '2.494.91'
$ws.Range("E18").Value = '  -0.31%  '

$ws.Range("D19").Value = '10.98'
$ws.Range("E19").Value = '  -2.27%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '7.18'
$ws.Range("E20").Value = '  -4.27%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '342.78'
$ws.Range("E21").Value = '  -1.63%  '

$ws.Range("D22").Value = '3.80'
$ws.Range("E22").Value = '  -2.84%  '

$ws.Range("E23").Value = '  +0.16%  '

$ws.Range("E24").Value = '  -3.51%  '

$ws.Range("D25").Value = '67.90'
$ws.Range("E25").Value = '  -3.67%  '

$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").Value = '1.11'
$ws.Range("E26").Value = '  +10.55%  '

$ws.Range("B27").Value = 'NEARProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D27").Value = '3.72'
$ws.Range("E27").Value = '  -5.98%  '

$ws.Range("E28").Value = '  -1.48%  '

$ws.Range("D29").Value = '8.18'
$ws.Range("E29").Value = '  -6.77%  '

$ws.Range("E30").Value = '  -6.24%  '

$ws.Range("D31").Value = '7.27'
$ws.Range("E31").Value = '  -7.12%  '

$ws.Range("D32").Value = '3.31'
$ws.Range("E32").Value = '  +123.59%  '

$ws.Range("E33").Value = '  -3.08%  '

$ws.Range("D34").Value = '433.71'
$ws.Range("E34").Value = '  -5.33%  '

$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("E36").Value = '  -3.25%  '

$ws.Range("D37").Value = '157.48'
$ws.Range("E37").Value = '  +0.36%  '

$ws.Range("D38").Value = '19.01'
$ws.Range("E38").Value = '  -0.28%  '

$ws.Range("E39").Value = '  -5.42%  '

$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("D41").Value = '17.89'
$ws.Range("E41").Value = '  -3.04%  '

$ws.Range("D42").Value = '0.307'
$ws.Range("E42").Value = '  -3.25%  '

$ws.Range("D43").Value = '4.47'
$ws.Range("E43").Value = '  -4.71%  '

$ws.Range("E44").Value = '  -5.01%  '

$ws.Range("D45").Value = '1.09'
$ws.Range("E45").Value = '  +0.46%  '

$ws.Range("E46").Value = '  -6.22%  '

$ws.Range("D47").Value = '134.24'
$ws.Range("E47").Value = '  -5.42%  '

$ws.Range("D48").Value = '3.36'
$ws.Range("E48").Value = '  -3.34%  '

$ws.Range("E49").Value = '  -2.06%  '

$ws.Range("D50").Value = '0.485'
$ws.Range("E50").Value = '  -6.65%  '

$ws.Range("E51").Value = '  -2.84%  '
